$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.710.16'
$ws.Range("E2").Value = '  -1.16%  '
$ws.Range("D3").Value = '3.778.47'
$ws.Range("E3").Value = '  +0.50%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '596.23'
$ws.Range("E5").Value = '  +0.46%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '167.07'
$ws.Range("E6").Value = '  -0.02%  '
$ws.Range("D7").Value = '3.777.02'
$ws.Range("E7").Value = '  +0.48%  '
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("E9").Value = '  -0.46%  '
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("E11").Value = '  -1.77%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.448'
$ws.Range("E12").Value = '  -0.13%  '
$ws.Range("E13").Value = '  -2.75%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '36.03'
$ws.Range("E14").Value = '  -0.28%  '
$ws.Range("D15").Value = '4.410.97'
$ws.Range("E15").Value = '  +0.45%  '
$ws.Range("D16").Value = '3.756.07'
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '67.662.14'
$ws.Range("E17").Value = '  -1.18%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '18.39'
$ws.Range("E18").Value = '  +2.23%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.01'
$ws.Range("E19").Value = '  +0.26%  '
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.04'
$ws.Range("E21").Value = '  -6.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '457.74'
$ws.Range("E22").Value = '  -1.56%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.695'
$ws.Range("E23").Value = '  -0.41%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.0000152'
$ws.Range("E24").Value = '  +2.23%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '83.23'
$ws.Range("E25").Value = '  -1.04%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.99'
$ws.Range("E26").Value = '  +0.17%  '
$ws.Range("E27").Value = '  -2.32%  '
$ws.Range("B28").Value = 'RenderToken'
$ws.Range("C28").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.03'
$ws.Range("E28").Value = '  -0.24%  '
$ws.Range("B29").Value = 'Dai'
$ws.Range("C29").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.77'
$ws.Range("E30").Value = '  -0.35%  '
$ws.Range("E31").Value = '  +3.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.22'
$ws.Range("E32").Value = '  -1.68%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '29.66'
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '9.11'
$ws.Range("E34").Value = '  -0.67%  '
$ws.Range("B35").Value = 'Binance-PegBSC-USD'
$ws.Range("C35").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.19%  '
$ws.Range("D36").Value = '3.729.87'
$ws.Range("E36").Value = '  +0.34%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.1000'
$ws.Range("E37").Value = '  -0.54%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.34'
$ws.Range("E38").Value = '  -1.33%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.138'
$ws.Range("E39").Value = '  +0.02%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.996'
$ws.Range("E40").Value = '  -0.35%  '
$ws.Range("E41").Value = '  -0.75%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '45.39'
$ws.Range("E44").Value = '  +3.22%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '48.16'
$ws.Range("E45").Value = '  +3.03%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.299'
$ws.Range("E46").Value = '  -1.28%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '149.47'
$ws.Range("E47").Value = '  +2.72%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.32'
$ws.Range("E48").Value = '  -2.01%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '390.62'
$ws.Range("E49").Value = '  +0.25%  '
$ws.Range("E50").Value = '  -4.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '26.14'
$ws.Range("E51").Value = '  +1.34%  '
